$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103; existing rows 103:121 shift down to 104:122
$ws.Rows("103:103").Insert()

# Populate the newly inserted row 103 with the new weekly entry
$ws.Range("A103").Value = 9
$ws.Range("B103").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C103").Value = "Metropolitana"
$ws.Range("D103").Value = 44504
$ws.Range("E103").Value = 13
$ws.Range("F103").Value = 100112003
$ws.Range("G103").Value = "Ajo"
$ws.Range("H103").Value = "Rosado"
$ws.Range("I103").Value = "1a nueva(o)"
$ws.Range("J103").Value = 36
$ws.Range("K103").Value = 3200
$ws.Range("L103").Value = 3400
$ws.Range("M103").Value = 3300
$ws.Range("N103").Value = "$/paquete 20 unidades (volumen en unidades)"
$ws.Range("O103").Value = "Provincia de Talagante"
$ws.Range("P103").Value = 165
$ws.Range("Q103").Value = 20
$ws.Range("R103").Value = "Hortaliza"
